$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 68

$ws.Cells.Item($row, 1).Value = "86KYC8"
$ws.Cells.Item($row, 2).Value = "Film de fusor Ricoh"
$ws.Cells.Item($row, 3).Value = "MPC 2000 2500 2800 3000 3300 3500 4500, SP810"
$ws.Cells.Item($row, 4).Value = 375000
$ws.Cells.Item($row, 5).Value = 650000
$ws.Cells.Item($row, 6).Value = 1
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Formula = "=(E68-D68)*G68"
$ws.Cells.Item($row, 9).Formula = "=D68*F68"
$ws.Cells.Item($row, 10).Value = 375000
